# Insert 2 new weekly data rows (Coliflor, Feria Lagunitas de Puerto Montt)
# right after the existing row 240, pushing the old rows 241-340 down to 243-342.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 241 (shifts old row 241.. down by 2).
$ws.Range("A241:A242").EntireRow.Insert()

# --- New row 241: Coliflor, Primera, Región Metropolitana, 2022-06-14 ---
$ws.Cells.Item(241, 1).Value = 4
$ws.Cells.Item(241, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(241, 3).Value = "Los Lagos"
$ws.Cells.Item(241, 4).Value = 44726
$ws.Cells.Item(241, 5).Value = 10
$ws.Cells.Item(241, 6).Value = 100112008
$ws.Cells.Item(241, 7).Value = "Coliflor"
$ws.Cells.Item(241, 8).Value = "Sin especificar"
$ws.Cells.Item(241, 9).Value = "Primera"
$ws.Cells.Item(241, 10).Value = 600
$ws.Cells.Item(241, 11).Value = 1800
$ws.Cells.Item(241, 12).Value = 1800
$ws.Cells.Item(241, 13).Value = 1800
$ws.Cells.Item(241, 14).Value = "`$/unidad"
$ws.Cells.Item(241, 15).Value = "Región Metropolitana"
$ws.Cells.Item(241, 16).Value = 1800
$ws.Cells.Item(241, 17).Value = 1
$ws.Cells.Item(241, 18).Value = "Hortaliza"

# --- New row 242: Coliflor, Segunda, Región Metropolitana, 2022-06-14 ---
$ws.Cells.Item(242, 1).Value = 4
$ws.Cells.Item(242, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(242, 3).Value = "Los Lagos"
$ws.Cells.Item(242, 4).Value = 44726
$ws.Cells.Item(242, 5).Value = 10
$ws.Cells.Item(242, 6).Value = 100112008
$ws.Cells.Item(242, 7).Value = "Coliflor"
$ws.Cells.Item(242, 8).Value = "Sin especificar"
$ws.Cells.Item(242, 9).Value = "Segunda"
$ws.Cells.Item(242, 10).Value = 600
$ws.Cells.Item(242, 11).Value = 1500
$ws.Cells.Item(242, 12).Value = 1500
$ws.Cells.Item(242, 13).Value = 1500
$ws.Cells.Item(242, 14).Value = "`$/unidad"
$ws.Cells.Item(242, 15).Value = "Región Metropolitana"
$ws.Cells.Item(242, 16).Value = 1500
$ws.Cells.Item(242, 17).Value = 1
$ws.Cells.Item(242, 18).Value = "Hortaliza"
